# Applies "Updated Results with corrected code" changes to Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D3: clear numeric value -> empty string cell (keep cell present, not fully blank)
$ws.Range("D3").Value = "'"
$ws.Range("D3").Style = "Normal"

# C4: 16042.61627423952 -> 0
$ws.Range("C4").Value = 0

# C5: 65484.51781727106 -> 0
$ws.Range("C5").Value = 0

# Row 7: rename "Other" -> "Biogas", update D7 value
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 66.31155363407012

# New row 8: "Other" with D8 value, copying formatting of row 7's A cell (bold/border/centered)
$ws.Range("A8").Value = "Other"
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)

$ws.Range("B8").Value = "'"
$ws.Range("B8").Style = "Normal"

$ws.Range("C8").Value = "'"
$ws.Range("C8").Style = "Normal"

$ws.Range("D8").Value = 47.52115294032427
